# ---------------------------------------------------------------------
# Applies the four changes described by the diff to
# "docs/Nuevas Funcionalidades.docx":
#
#  1. Remove the paragraph "TERMINADO (SOLO ME FALTA EDITAR DETALLE ...)"
#  2. Add a <w:lastRenderedPageBreak/> marker before "Transferencia"
#  3. Change the "Editar detalle" paragraph shading from FFD966 (themed)
#     to a plain C00000 fill
#  4. Replace the trailing empty paragraph with three new bullet items
#     ("Faltan varios informes", "Editar transferencia",
#     "Mensajes/notificaciones(mejorarlas)") plus a blank spacer line
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-PkgXml([string]$bodyInnerXml) {
    return '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParagraphIndex([string]$literalPrefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like ($literalPrefix + "*")) {
            return $i
        }
    }
    return -1
}

# --- 1. Delete the "TERMINADO (...)" paragraph entirely -----------------
$idxTerminado = Find-ParagraphIndex "TERMINADO (SOLO ME FALTA EDITAR DETALLE"
if ($idxTerminado -gt 0) {
    $d.Paragraphs.Item($idxTerminado).Range.Delete()
}

# --- 2. Add <w:lastRenderedPageBreak/> before "Transferencia" -----------
$idxTransferencia = Find-ParagraphIndex "Transferencia"
if ($idxTransferencia -gt 0) {
    $pTrans = $d.Paragraphs.Item($idxTransferencia)
    $rTrans = $pTrans.Range.Duplicate
    $rTrans.Collapse(1)
    $inner = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Transferencia</w:t></w:r></w:p>'
    $rTrans.InsertXML((Get-PkgXml $inner))
}

# --- 3. Change "Editar detalle" paragraph shading to plain C00000 -------
$idxEditarDetalle = Find-ParagraphIndex "Editar detalle"
if ($idxEditarDetalle -gt 0) {
    $pEd = $d.Paragraphs.Item($idxEditarDetalle)
    $rEd = $pEd.Range.Duplicate
    $rEd.Collapse(1)
    $inner = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
             '<w:shd w:val="clear" w:color="auto" w:fill="C00000"/>' +
             '<w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Editar detalle</w:t></w:r>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:tab/></w:r>' +
             '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> X</w:t></w:r></w:p>'
    $rEd.InsertXML((Get-PkgXml $inner))
}

# --- 4. Replace the trailing empty paragraph with the three new items ---
#        plus the whitespace-only spacer paragraph.
$idxPruebas = Find-ParagraphIndex "Pruebas"
$idxTarget = $idxPruebas + 1
$pTarget = $d.Paragraphs.Item($idxTarget)
$rTarget = $pTarget.Range.Duplicate
$rTarget.Collapse(1)

$p1 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="C00000"/>' +
      '<w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr>' +
      '<w:t>Faltan varios informes                                                                              X</w:t></w:r></w:p>'

$p2 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="C00000"/>' +
      '<w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr>' +
      '<w:t>Editar transferencia                                                                                  X</w:t></w:r></w:p>'

$p3 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
      '<w:shd w:val="clear" w:color="auto" w:fill="C00000"/>' +
      '<w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>Mensajes/notificaciones(</w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">mejorarlas)   </w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">                                                 X</w:t></w:r></w:p>'

$p4 = '<w:p><w:pPr><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve">                              </w:t></w:r></w:p>'

$inner = $p1 + $p2 + $p3 + $p4
$rTarget.InsertXML((Get-PkgXml $inner))
